# WIP_13.2.pptx — title slide text update
#
# The deck's title slide ("Title 1" placeholder on slide 1) is renamed from
# "Kinetics in Fixed Axis Rotation Systems" to
# "Kinetics in Belt and Gear Driven Systems" to reflect the corrected topic
# for this lecture in the renumbered lecture sequence.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Kinetics in Belt and Gear Driven Systems"
